$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.088846
$ws.Range("H2").Value = 30.266538
$ws.Range("I2").Value = 0.1151445838515654
$ws.Range("J2").Value = 0.1151445838515654
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 95.07579491774068
$ws.Range("R2").Value = 855.682154259666
$ws.Range("S2").Value = 0.007662107829449277
$ws.Range("T2").Value = 0.007662107829449277
$ws.Range("G3").Value = 10.088846
$ws.Range("H3").Value = 30.266538
$ws.Range("I3").Value = 0.1151445838515654
$ws.Range("J3").Value = 0.1151445838515654
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 510.4205936670414
$ws.Range("R3").Value = 4593.785343003372
$ws.Range("S3").Value = 0.04113452462250864
$ws.Range("T3").Value = 0.04113452462250864
$ws.Range("G4").Value = 10.088846
$ws.Range("H4").Value = 30.266538
$ws.Range("I4").Value = 0.1151445838515654
$ws.Range("J4").Value = 0.1151445838515654
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 270.7923210665493
$ws.Range("R4").Value = 2437.130889598944
$ws.Range("S4").Value = 0.02182300937051219
$ws.Range("T4").Value = 0.02182300937051219
$ws.Range("G5").Value = 10.088846
$ws.Range("H5").Value = 30.266538
$ws.Range("I5").Value = 0.1151445838515654
$ws.Range("J5").Value = 0.1151445838515654
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 552.4908225399939
$ws.Range("R5").Value = 4972.417402859945
$ws.Range("S5").Value = 0.04452494202909531
$ws.Range("T5").Value = 0.04452494202909531
$ws.Range("I6").Value = 0.4327250566572728
$ws.Range("J6").Value = 0.4327250566572729
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 357.3045067890552
$ws.Range("R6").Value = 3215.740561101496
$ws.Range("S6").Value = 0.028794980481989
$ws.Range("T6").Value = 0.02879498048198901
$ws.Range("I7").Value = 0.4327250566572728
$ws.Range("J7").Value = 0.4327250566572729
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("S7").Value = 0.1545877270336154
$ws.Range("T7").Value = 0.1545877270336154
$ws.Range("I8").Value = 0.4327250566572728
$ws.Range("J8").Value = 0.4327250566572729
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 1017.665082944185
$ws.Range("R8").Value = 9158.985746497665
$ws.Range("S8").Value = 0.08201308867867081
$ws.Range("T8").Value = 0.08201308867867084
$ws.Range("I9").Value = 0.4327250566572728
$ws.Range("J9").Value = 0.4327250566572729
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 2076.316701048131
$ws.Range("R9").Value = 18686.85030943318
$ws.Range("S9").Value = 0.1673292604629976
$ws.Range("T9").Value = 0.1673292604629976
$ws.Range("G10").Value = 15.69885766666667
$ws.Range("H10").Value = 47.096573
$ws.Range("I10").Value = 0.1791719719949428
$ws.Range("J10").Value = 0.1791719719949428
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 147.9437164526846
$ws.Range("R10").Value = 1331.493448074161
$ws.Range("S10").Value = 0.01192270555434947
$ws.Range("T10").Value = 0.01192270555434947
$ws.Range("G11").Value = 15.69885766666667
$ws.Range("H11").Value = 47.096573
$ws.Range("I11").Value = 0.1791719719949428
$ws.Range("J11").Value = 0.1791719719949428
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 794.2454716936292
$ws.Range("R11").Value = 7148.209245242661
$ws.Range("S11").Value = 0.06400782083845452
$ws.Range("T11").Value = 0.06400782083845452
$ws.Range("G12").Value = 15.69885766666667
$ws.Range("H12").Value = 47.096573
$ws.Range("I12").Value = 0.1791719719949428
$ws.Range("J12").Value = 0.1791719719949428
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 421.3693127687804
$ws.Range("R12").Value = 3792.323814919024
$ws.Range("S12").Value = 0.03395792917901649
$ws.Range("T12").Value = 0.03395792917901649
$ws.Range("G13").Value = 15.69885766666667
$ws.Range("H13").Value = 47.096573
$ws.Range("I13").Value = 0.1791719719949428
$ws.Range("J13").Value = 0.1791719719949428
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 859.7093052262823
$ws.Range("R13").Value = 7737.38374703654
$ws.Range("S13").Value = 0.06928351642312232
$ws.Range("T13").Value = 0.06928351642312232
$ws.Range("G14").Value = 23.91632366666667
$ws.Range("H14").Value = 71.748971
$ws.Range("I14").Value = 0.2729583874962189
$ws.Range("J14").Value = 0.2729583874962189
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 225.3839025908719
$ws.Range("R14").Value = 2028.455123317847
$ws.Range("S14").Value = 0.0181635690363407
$ws.Range("T14").Value = 0.0181635690363407
$ws.Range("G15").Value = 23.91632366666667
$ws.Range("H15").Value = 71.748971
$ws.Range("I15").Value = 0.2729583874962189
$ws.Range("J15").Value = 0.2729583874962189
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 1209.988151694764
$ws.Range("R15").Value = 10889.89336525287
$ws.Range("S15").Value = 0.09751230267033376
$ws.Range("T15").Value = 0.09751230267033376
$ws.Range("G16").Value = 23.91632366666667
$ws.Range("H16").Value = 71.748971
$ws.Range("I16").Value = 0.2729583874962189
$ws.Range("J16").Value = 0.2729583874962189
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 641.9323673961832
$ws.Range("R16").Value = 5777.391306565648
$ws.Range("S16").Value = 0.05173298863773607
$ws.Range("T16").Value = 0.05173298863773608
$ws.Range("G17").Value = 23.91632366666667
$ws.Range("H17").Value = 71.748971
$ws.Range("I17").Value = 0.2729583874962189
$ws.Range("J17").Value = 0.2729583874962189
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 1309.718607532456
$ws.Range("R17").Value = 11787.4674677921
$ws.Range("S17").Value = 0.1055495271518084
$ws.Range("T17").Value = 0.1055495271518084
